# Update monthly Rugvista trend values in column B (re-run with local data, 20 Oct 2025)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 15
$ws.Range("B10").Value = 14
$ws.Range("B27").Value = 21
$ws.Range("B32").Value = 16
$ws.Range("B33").Value = 17
$ws.Range("B44").Value = 17
$ws.Range("B49").Value = 24
$ws.Range("B51").Value = 28
$ws.Range("B52").Value = 22
$ws.Range("B53").Value = 27
$ws.Range("B55").Value = 25
$ws.Range("B57").Value = 29
$ws.Range("B64").Value = 43
$ws.Range("B66").Value = 29
$ws.Range("B69").Value = 28
$ws.Range("B72").Value = 43
$ws.Range("B75").Value = 38
$ws.Range("B86").Value = 43
$ws.Range("B89").Value = 33
$ws.Range("B90").Value = 29
$ws.Range("B91").Value = 25
$ws.Range("B92").Value = 32
$ws.Range("B97").Value = 50
$ws.Range("B98").Value = 57
$ws.Range("B100").Value = 47
$ws.Range("B104").Value = 38
$ws.Range("B106").Value = 55
$ws.Range("B107").Value = 76
$ws.Range("B109").Value = 72
$ws.Range("B110").Value = 84
$ws.Range("B112").Value = 73
$ws.Range("B119").Value = 86
